$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# Lift Coefficient CL: -2 -> -1.98
$ws.Range("C8").Value = -1.98

# Drag Coefficient CD: -1.2 -> -1.33
$ws.Range("C9").Value = -1.33

# Front Aero Distribution: 47 -> formula =100-56.3 (result 43.7)
$ws.Range("C12").Formula = "=100-56.3"

# Frontal Area: 1.1 -> 1.15
$ws.Range("C13").Value = 1.15

# Update the selected cell on the Info sheet to E8
$ws.Activate()
$ws.Range("E8").Select()
